$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance")
$ws.Range("A1").Value = "TEST"
